# Correct the French word on the "Feuil1" sheet: "air" -> "eau"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Range("A1").Value = "eau"
